# Timesheet January 2020 - add timesheet entries for 26-02-2020
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- First lay down formatting for all the new rows (copy from matching existing rows) ----

# Row 209: blank separator row, like the other day-boundary rows (e.g. row 198)
$ws.Range("A198:C198").Copy()
$ws.Range("A209:C209").PasteSpecial(-4122)

# Rows 210-211: wrapped two-line task rows, like row 19
$ws.Range("A19:C19").Copy()
$ws.Range("A210:C211").PasteSpecial(-4122)
$ws.Rows("210:211").RowHeight = 30

# Rows 212-213: plain single-line rows, like row 10
$ws.Range("A10:C10").Copy()
$ws.Range("A212:C213").PasteSpecial(-4122)

# Row 214: lunch row, formatted like the other lunch row that uses the date style on column A (row 54)
$ws.Range("A54:C54").Copy()
$ws.Range("A214:C214").PasteSpecial(-4122)

# Rows 215-218: plain single-line rows, like row 10
$ws.Range("A10:C10").Copy()
$ws.Range("A215:C218").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Now fill in the values, in the order they were actually typed ----

$ws.Range("B210").Value = "Performing multi variate timeseries on combined horly data, done upto making date" + [char]10 + "time index."
$ws.Range("B211").Value = "Build Multivariate timeseries model using var algo. Got problem in label decoding." + [char]10 + "Results were not good. Working on sample UI page."
$ws.Range("B212").Value = "made wireframe of sample UI"
$ws.Range("B213").Value = "viewd some django documentation"

$ws.Range("A210").Value = "Feb 26 10:00 to 11:00"
$ws.Range("A211").Value = "Feb 26 11:00 to 12:00"
$ws.Range("A212").Value = "Feb 26 12:00 to 13:00"
$ws.Range("A213").Value = "Feb 26 13:00 to 13:30"
$ws.Range("A214").Value = "Feb 26 13:30 to 14:00"
$ws.Range("B214").Value = "Lunch"

$ws.Range("A215").Value = "Feb 26 14:00 to 15:00"
$ws.Range("B215").Value = "Working on rest api"

$ws.Range("A216").Value = "Feb 26 15:00 to 16:00"
$ws.Range("B216").Value = "Made a sample django api with rest, used post for json sucessfully"

$ws.Range("A217").Value = "Feb 26 16:00 to 17:00"
$ws.Range("B217").Value = "Made sample django app which accepts initials of name and return full name"

$ws.Range("A218").Value = "Feb 26 17:00 to 20:00"
$ws.Range("B218").Value = "Did a example considering up time, output count and harmful alarms"

# Column C ("Infimetrics") for every new data row
$ws.Range("C210").Value = "Infimetrics"
$ws.Range("C211").Value = "Infimetrics"
$ws.Range("C212").Value = "Infimetrics"
$ws.Range("C213").Value = "Infimetrics"
$ws.Range("C214").Value = "Infimetrics"
$ws.Range("C215").Value = "Infimetrics"
$ws.Range("C216").Value = "Infimetrics"
$ws.Range("C217").Value = "Infimetrics"
$ws.Range("C218").Value = "Infimetrics"

# ---- Update the view: scrolled so row 195 is top, D218 selected ----
$ws.Range("D218").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 195 | Out-Null
